# Regenerate the handoff/handback report: a new GUID-named source file
# ("4cc78c11-207a-45da-a586-a3d051d1b02c.md") replaces the old one
# ("80317823-dd0e-4c0d-9dc1-b9aaaeaa6b0e.md"), the xliff handoff files get
# new content hashes, and the associated timestamps move forward.

$wb = $excel.ActiveWorkbook

$newGuid = "4cc78c11-207a-45da-a586-a3d051d1b02c"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-10-21 00:50:39"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.7fdb82234c0b54e8c728030112558f41a17f29d8.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-10-21 00:50:27"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.7fdb82234c0b54e8c728030112558f41a17f29d8.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-10-21 00:50:39"
